$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008 and 2009 data rows (rows 2 and 3); everything below shifts up two rows.
$ws.Rows("2:3").Delete()

# Append the new 2021 data row (now row 13, right after 2020 in row 12).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13").Value = "2021年"

$ws.Range("B13").Value = 12208
$ws.Range("C13").Value = 79360
$ws.Range("D13").Value = 116915
$ws.Range("E13").Value = 411340
$ws.Range("F13").Value = 240569
$ws.Range("G13").Value = 95280
$ws.Range("H13").Value = 520053
$ws.Range("I13").Value = 29353
$ws.Range("L13").Value = 403138
$ws.Range("M13").Value = 171996
